$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "site" sheet at the front of the workbook.
# ---------------------------------------------------------------------------
$site = $wb.Worksheets.Add()
$site.Name = "site"
$site.Move($wb.Worksheets.Item(1))

# ---------------------------------------------------------------------------
# 2. Make room for the new "site_name" / "guide" join columns.
# ---------------------------------------------------------------------------
$outing = $wb.Worksheets.Item("outing")
$outing.Columns.Item(3).Insert()

$capture = $wb.Worksheets.Item("capture")
$capture.Columns.Item(3).Insert()

# ---------------------------------------------------------------------------
# 3. Populate cell data (kept in the same order the shared-string table in
#    the target file implies: outing!C7 ("site") first, then the site sheet /
#    outing!C column ("site_name", "Bendy Bay", "name of the site").
# ---------------------------------------------------------------------------
$outing.Range("C7").Value = "site"

$site.Range("A1").Value = "name"
$site.Range("B1").Value = "site_name"
$site.Range("A2").Value = "example"
$site.Range("B2").Value = "Bendy Bay"
$site.Range("A3").Value = "description"
$site.Range("B3").Value = "name of the site"
$site.Range("A4").Value = "chk"
$site.Range("B4").Value = 'c("")'
$site.Range("A5").Value = "pkey"
$site.Range("B5").Value = $true
$site.Range("A6").Value = "unique"
$site.Range("A7").Value = "join"

$outing.Range("C1").Value = "site_name"
$outing.Range("C2").Value = "Bendy Bay"
$outing.Range("C3").Value = "name of the site"
$outing.Range("C4").Value = 'c("")'

# ---------------------------------------------------------------------------
# 4. Populate the new "guide" column on the capture sheet (all reuse existing
#    shared strings already present in the workbook).
# ---------------------------------------------------------------------------
$capture.Range("C1").Value = "guide"
$capture.Range("C2").Value = "GA"
$capture.Range("C3").Value = "initials of guide"
$capture.Range("C4").Value = 'c("GA")'
$capture.Range("C5").Value = $true
$capture.Range("C7").Value = "outing"

# ---------------------------------------------------------------------------
# 5. Fix up the pkey / unique marker flags that moved columns.
# ---------------------------------------------------------------------------
$outing.Range("I5").Value = $true
$outing.Range("B6").Value = ""
$outing.Range("F6").Value = $true

# ---------------------------------------------------------------------------
# 5b. Column widths for the two freshly inserted columns.
# ---------------------------------------------------------------------------
$outing.Columns.Item(3).ColumnWidth = 15.75
$capture.Columns.Item(3).ColumnWidth = 19.084

# ---------------------------------------------------------------------------
# 6. Selections per sheet.
# ---------------------------------------------------------------------------
$site.Range("B6").Select()
$outing.Columns.Item(9).Select()

$recapture = $wb.Worksheets.Item("recapture")
$recapture.Range("H5").Select()

$capture.Range("C8").Select()
$capture.Activate()
